$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The diff updates every data row's "Förändrad" (column C) date from
# 45188 (2023-09-19) to 45189 (2023-09-20), across rows 2 through 453.
$ws.Range("C2:C453").Value = 45189
